# Update the "想去人数" (want-to-go count) values in column F for the
# rows that changed between the two scraped snapshots. The same set of
# updates applies identically to the "展览" and "全部类型" worksheets,
# which carry duplicate data in this workbook.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value (applies to both affected sheets)
$updates = @{
    "F3"  = 3099
    "F5"  = 2633
    "F9"  = 1371
    "F13" = 1183
    "F14" = 356
    "F15" = 326
    "F16" = 37
    "F21" = 2499
    "F22" = 30
    "F23" = 281
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
